$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full, final ordered list of gene names for column A (rows 1..162).
# Rows 1-141 keep their original row numbers but some now show different
# genes because new unique gene names were inserted into the shared-string
# table ahead of them (alphabetically, within each existing block).
# Rows 142-162 are newly appended rows that repeat the (now shifted)
# trailing gene names.
$values = @(
    "ACP5",
    "ACVR1B",
    "B4GALT1",
    "BAIAP2",
    "BAX",
    "CA2",
    "CCL4",
    "CD8A",
    "CD8B",
    "CDKAL1",
    "CDKN1C",
    "COX15",
    "CTSB",
    "ELOVL6",
    "FSTL1",
    "GM2A",
    "GYPB",
    "HIST1H1C",
    "HLA-DOB",
    "IGF1R",
    "JUP",
    "KLF1",
    "MAK",
    "MSC",
    "OLFM1",
    "PGM1",
    "PLA2G16",
    "SDHD",
    "SMC2",
    "ST14",
    "SYNJ2",
    "TIPIN",
    "TMEM176B",
    "TNNT1",
    "TUBB2A",
    "UQCRB",
    "FCER1A",
    "MYOM2",
    "S100B",
    "ABCB4",
    "CD22",
    "COCH",
    "FADS2",
    "FGFR1",
    "GLTP",
    "GSTT1",
    "HAMP",
    "IGLC1",
    "MBP",
    "NKX3-1",
    "PCGF3",
    "PLVAP",
    "QRSL1",
    "RHD",
    "RPL37A",
    "SPON1",
    "TCL1A",
    "TPPP3",
    "VPREB3",
    "WASL",
    "ADIPOR2",
    "ATP6V0E1",
    "ATRX",
    "BCL2A1",
    "BGN",
    "CD40",
    "CDCA8",
    "CRAT",
    "DDX39A",
    "E2F2",
    "EIF4H",
    "ERBB2",
    "ESR1",
    "EZR",
    "F7",
    "FSTL3",
    "GLIPR1",
    "GPX1",
    "GRHPR",
    "GYPA",
    "GYPE",
    "HGF",
    "HIST3H2A",
    "HNF4A",
    "IDH3B",
    "IGF1",
    "IL16",
    "LMNA",
    "LTB",
    "MARCH6",
    "MCM4",
    "MYH11",
    "NDUFB7",
    "NEK2",
    "NF1",
    "NR3C1",
    "PA2G4",
    "PDLIM5",
    "PLEK",
    "PRMT2",
    "PRPF31",
    "RPS19",
    "RRAD",
    "SFPQ",
    "SMAD3",
    "SORBS3",
    "TPD52L1",
    "UBE2L3",
    "ATG2A",
    "CAMTA2",
    "DIP2A",
    "HYOU1",
    "NSUN6",
    "OTUD7B",
    "POLD3",
    "PSTPIP1",
    "RRNAD1",
    "SCML1",
    "SMARCD2",
    "TCIRG1",
    "ANK1",
    "ARHGEF12",
    "BUB3",
    "CANX",
    "CAPZB",
    "CASP1",
    "CD36",
    "CD74",
    "CDK2",
    "CELF2",
    "ELN",
    "EPB42",
    "F2RL1",
    "FBXO9",
    "FLNA",
    "HS2ST1",
    "KYNU",
    "MARCKS",
    "NCBP1",
    "NCOA4",
    "PCLO",
    "PRF1",
    "PRPF4B",
    "RBM5",
    "SLC2A3",
    "SLC4A1",
    "SLC7A5",
    "SPTB",
    "SULT1A1",
    "TGM2",
    "TRIM58",
    "UBE2D1",
    "ZFP36L1",
    "DCAF8",
    "DYNLL1",
    "HIST1H1T",
    "HLA-B",
    "IL2RG",
    "MUC8",
    "NCOA2",
    "PHLPP1",
    "RUBCNL"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

Write-Output ("Rows written: " + $values.Length)
Write-Output ("A1=" + $ws.Cells.Item(1,1).Value2)
Write-Output ("A142=" + $ws.Cells.Item(142,1).Value2)
Write-Output ("A162=" + $ws.Cells.Item(162,1).Value2)
